$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2978.1428
$ws.Range("J17").Value = 2978.1428
$ws.Range("L17").Value = 8934.428400000001
$ws.Range("N17").Value = -9270.428400000001

$ws.Range("H38").Value = 2203.5
$ws.Range("I38").Value = 1103.1666
$ws.Range("J38").Value = 5504.5
$ws.Range("K38").Value = 3309.4998
$ws.Range("L38").Value = 16513.5
$ws.Range("M38").Value = -2937.4998
$ws.Range("N38").Value = -17257.5

$ws.Range("H39").Value = 63.5
$ws.Range("I39").Value = 36.22222
$ws.Range("J39").Value = 309
$ws.Range("K39").Value = 108.66666
$ws.Range("L39").Value = 927
$ws.Range("M39").Value = 187.33334
$ws.Range("N39").Value = -1519

$ws.Range("H40").Value = 6961.353
$ws.Range("J40").Value = 8839.4
$ws.Range("L40").Value = 8839.4
$ws.Range("N40").Value = -9189.4

$ws.Range("H51").Value = 5500
$ws.Range("J51").Value = 5500
$ws.Range("L51").Value = 5500
$ws.Range("N51").Value = -6468

$ws.Range("H92").Value = 4808761.5
$ws.Range("I92").Value = 789.3333
$ws.Range("K92").Value = 789.3333
$ws.Range("M92").Value = 458.6667

$ws.Range("H101").Value = 1977.2222
$ws.Range("I101").Value = 1218.3334
$ws.Range("K101").Value = 3655.0002
$ws.Range("M101").Value = -2033.0002

$ws.Range("H112").Value = 2150.6072
$ws.Range("J112").Value = 2212.72
$ws.Range("L112").Value = 6638.16
$ws.Range("N112").Value = -8854.16

$ws.Range("H137").Value = 1993.5405
$ws.Range("I137").Value = 1524.3182
$ws.Range("K137").Value = 4572.9546
$ws.Range("M137").Value = -2022.9546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14357.909
$ws.Range("I2").Value = 19364.875
$ws.Range("K2").Value = 19364.875
$ws.Range("M2").Value = -19251.875

$ws.Range("H32").Value = 15636287
$ws.Range("I32").Value = 16675423
$ws.Range("K32").Value = 16675423
$ws.Range("M32").Value = -16675136

$ws.Range("H45").Value = 4749.5713
$ws.Range("I45").Value = 3236.6667
$ws.Range("K45").Value = 3236.6667
$ws.Range("M45").Value = -2859.6667

$ws.Range("H61").Value = 3735.1702
$ws.Range("I61").Value = 2930.8286
$ws.Range("J61").Value = 6081.1665
$ws.Range("K61").Value = 2930.8286
$ws.Range("L61").Value = 6081.1665
$ws.Range("M61").Value = -2718.8286
$ws.Range("N61").Value = -6505.1665

$ws.Range("H116").Value = 14357.909
$ws.Range("I116").Value = 19364.875
$ws.Range("K116").Value = 19364.875
$ws.Range("M116").Value = -17070.875

$ws.Range("H136").Value = 3735.1702
$ws.Range("I136").Value = 2930.8286
$ws.Range("J136").Value = 6081.1665
$ws.Range("K136").Value = 8792.485799999999
$ws.Range("L136").Value = 18243.4995
$ws.Range("M136").Value = -6242.485799999999
$ws.Range("N136").Value = -23343.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14357.909
$ws.Range("I3").Value = 19364.875
$ws.Range("K3").Value = 19364.875
$ws.Range("M3").Value = -19250.875

$ws.Range("H105").Value = 4405.5884
$ws.Range("I105").Value = 2551.25
$ws.Range("K105").Value = 2551.25
$ws.Range("M105").Value = -804.25

$ws.Range("H134").Value = 2889.2432
$ws.Range("I134").Value = 1868.9231
$ws.Range("K134").Value = 5606.7693
$ws.Range("M134").Value = -3071.7693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1366.6666
$ws.Range("I6").Value = 2250.5
$ws.Range("J6").Value = 1114.1428
$ws.Range("K6").Value = 2250.5
$ws.Range("L6").Value = 1114.1428
$ws.Range("M6").Value = -2137.5
$ws.Range("N6").Value = -1340.1428

$ws.Range("H31").Value = 3160.878
$ws.Range("I31").Value = 2369.516
$ws.Range("J31").Value = 5614.1
$ws.Range("K31").Value = 2369.516
$ws.Range("L31").Value = 5614.1
$ws.Range("M31").Value = -2074.516
$ws.Range("N31").Value = -6204.1

$ws.Range("H34").Value = 3160.878
$ws.Range("I34").Value = 2369.516
$ws.Range("J34").Value = 5614.1
$ws.Range("K34").Value = 2369.516
$ws.Range("L34").Value = 5614.1
$ws.Range("M34").Value = -2167.516
$ws.Range("N34").Value = -6018.1

$ws.Range("H41").Value = 898.625
$ws.Range("I41").Value = 898.625
$ws.Range("K41").Value = 898.625
$ws.Range("M41").Value = -470.625

$ws.Range("H58").Value = 2107.2222
$ws.Range("I58").Value = 1103.1
$ws.Range("K58").Value = 1103.1
$ws.Range("M58").Value = -900.0999999999999

$ws.Range("H99").Value = 16034792
$ws.Range("J99").Value = 33339042
$ws.Range("L99").Value = 33339042
$ws.Range("N99").Value = -33342038

$ws.Range("H126").Value = 16034792
$ws.Range("J126").Value = 33339042
$ws.Range("L126").Value = 100017126
$ws.Range("N126").Value = -100022066

$ws.Range("H136").Value = 2107.2222
$ws.Range("I136").Value = 1103.1
$ws.Range("K136").Value = 3309.3
$ws.Range("M136").Value = -759.2999999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 163.3158
$ws.Range("I2").Value = 57.5
$ws.Range("K2").Value = 345
$ws.Range("M2").Value = -232

$ws.Range("H11").Value = 382.77777
$ws.Range("I11").Value = 349.2857
$ws.Range("K11").Value = 1047.8571
$ws.Range("M11").Value = -907.8571000000002

$ws.Range("H12").Value = 238.33333
$ws.Range("J12").Value = 257.23077
$ws.Range("L12").Value = 771.69231
$ws.Range("N12").Value = -1117.69231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5275.75
$ws.Range("I126").Value = 4151.5
$ws.Range("K126").Value = 12454.5
$ws.Range("M126").Value = -9984.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4263.4116
$ws.Range("I7").Value = 2427.875
$ws.Range("K7").Value = 2427.875
$ws.Range("M7").Value = -2315.875

$ws.Range("H40").Value = 8969.08
$ws.Range("I40").Value = 9915.857
$ws.Range("K40").Value = 9915.857
$ws.Range("M40").Value = -9779.857

$ws.Range("H61").Value = 4382.6665
$ws.Range("I61").Value = 2608.5
$ws.Range("K61").Value = 2608.5
$ws.Range("M61").Value = -2406.5

$ws.Range("H93").Value = 14446718
$ws.Range("I93").Value = 2070.7693
$ws.Range("J93").Value = 52002800
$ws.Range("K93").Value = 2070.7693
$ws.Range("L93").Value = 52002800
$ws.Range("M93").Value = -822.7692999999999
$ws.Range("N93").Value = -52005296

$ws.Range("H100").Value = 65318.39
$ws.Range("I100").Value = 112733.6
$ws.Range("K100").Value = 112733.6
$ws.Range("M100").Value = -112192.6

$ws.Range("H107").Value = 4470
$ws.Range("I107").Value = 4470
$ws.Range("K107").Value = 4470
$ws.Range("M107").Value = -2550

$ws.Range("H113").Value = 4382.6665
$ws.Range("I113").Value = 2608.5
$ws.Range("K113").Value = 2608.5
$ws.Range("M113").Value = -438.5

$ws.Range("H122").Value = 4739.0967
$ws.Range("I122").Value = 3671.8096
$ws.Range("J122").Value = 6980.4
$ws.Range("K122").Value = 11015.4288
$ws.Range("L122").Value = 20941.2
$ws.Range("M122").Value = -8565.4288
$ws.Range("N122").Value = -25841.2

$ws.Range("H126").Value = 4263.4116
$ws.Range("I126").Value = 2427.875
$ws.Range("K126").Value = 7283.625
$ws.Range("M126").Value = -4813.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10103990
$ws.Range("I81").Value = 2682.6191
$ws.Range("J81").Value = 27781278
$ws.Range("K81").Value = 5365.2382
$ws.Range("L81").Value = 55562556
$ws.Range("M81").Value = -4304.2382
$ws.Range("N81").Value = -55564678

$ws.Range("H84").Value = 10103990
$ws.Range("I84").Value = 2682.6191
$ws.Range("J84").Value = 27781278
$ws.Range("K84").Value = 26826.191
$ws.Range("L84").Value = 277812780
$ws.Range("M84").Value = -21522.191
$ws.Range("N84").Value = -277823388

$ws.Range("H96").Value = 33529.938
$ws.Range("I96").Value = 40659.617
$ws.Range("J96").Value = 2634.6667
$ws.Range("K96").Value = 40659.617
$ws.Range("L96").Value = 2634.6667
$ws.Range("M96").Value = -39286.617
$ws.Range("N96").Value = -5380.6667

$ws.Range("H125").Value = 69375
$ws.Range("J125").Value = 69375
$ws.Range("L125").Value = 69375
$ws.Range("N125").Value = -79215

$ws.Range("H136").Value = 1867.807
$ws.Range("I136").Value = 951
$ws.Range("J136").Value = 5701.727
$ws.Range("K136").Value = 2853
$ws.Range("L136").Value = 17105.181
$ws.Range("M136").Value = -303
$ws.Range("N136").Value = -22205.181
